$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a bare number-looking string (e.g. "241.54").
# The source keeps these as literal text cells (inline strings), but a
# plain Value= assignment would let Excel auto-coerce them into real
# numbers. Force Text format before writing, then restore the default
# "Normal" style afterwards so the cell serializes without an explicit
# style index, matching the rest of the sheet.
$textForceCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D16",
    "D19",
    "D20",
    "D23",
    "D24",
    "D25",
    "D26",
    "D28",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D42",
    "D46",
    "D47",
    "D48",
    "D50",
    "D51"
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Refreshed price (D) and volume-change (E) figures, plus row 51 which
# is rebranded from Cronos to Algorand (coin name, link, price, change).
$ws.Range("D2").Value = "30.021.60"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").Value = "1.905.64"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "0.7606"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").Value = "241.54"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "0.3082"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("D9").Value = "25.57"
$ws.Range("E9").Value = "  -6.62%  "
$ws.Range("D10").Value = "0.06904"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").Value = "0.08021"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "0.7554"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").Value = "1.903.14"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").Value = "5.253"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "91.83"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "6.191"
$ws.Range("E16").Value = "  +5.33%  "
$ws.Range("D17").Value = "30.026.81"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").Value = "0.000007753"
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("D20").Value = "237.65"
$ws.Range("E20").Value = "  -3.62%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").Value = "2.146.90"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "0.9996"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "7.041"
$ws.Range("E24").Value = "  +5.82%  "
$ws.Range("D25").Value = "9.318"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("D26").Value = "166.63"
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").Value = "0.1295"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("E29").Value = "  -3.41%  "
$ws.Range("D30").Value = "1.344"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "1.528"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").Value = "4.308"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("D33").Value = "4.046"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("E34").Value = "  +3.86%  "
$ws.Range("D35").Value = "1.288"
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("D36").Value = "0.7378"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("D38").Value = "0.01946"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "2.764"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").Value = "6.252"
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").Value = "72.80"
$ws.Range("E42").Value = "  -4.32%  "
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").Value = "7.707"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").Value = "101.63"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").Value = "9.889"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "2.051.49"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Value = "36.55"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.1162"
$ws.Range("E51").Value = "  -4.46%  "

# Restore default styling on the text-forced cells.
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
